# "Dokument 4. iterace plus přerozdělení bodů."
#
# The 4th-iteration work log sheet gets the previously empty H/I (Petr)
# and J/K (Lenka) cells filled in for the two existing "10.týden" rows,
# two brand new task rows are inserted before the old "11.týden" row, and
# a final grand-total row is appended that sums every iteration together.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List1")

# ---------------------------------------------------------------------
# 1) Fill in the previously empty H/I (Petr) and J/K (Lenka) cells for
#    the existing "10.týden" rows (35-37). Cell styles are already
#    correct here (16 / 32), so plain value assignment is enough.
# ---------------------------------------------------------------------
$ws.Range("H35").Value = 0.25
$ws.Range("I35").Value = "Oprava chyb v generování sekvenčních diagramů"

$ws.Range("H36").Value = 2
$ws.Range("I36").Value = "Oprava databázového modelu dle oponentury"
$ws.Range("J36").Value = 0.5
$ws.Range("K36").Value = "Sepsání zprávy o implementaci."

$ws.Range("H37").Value = 3
$ws.Range("I37").Value = "Hledání, jak správně vytvořit Návrh tříd. Pokus o vytvoření návrhu jedné části."
$ws.Range("J37").Value = 4
$ws.Range("K37").Value = "Oprava modelu nasazení"

# ---------------------------------------------------------------------
# 2) Insert two brand-new rows before the old row 38 ("11.týden"), which
#    pushes it (and the "Celkem ke 4. iteraci" row after it) down by two.
#    Copy the formatting of the row directly above down into both new
#    rows so the table borders/alignment stay consistent.
# ---------------------------------------------------------------------
$ws.Rows.Item(38).Insert()
$ws.Rows.Item(39).Insert()

$ws.Range("A37:K37").Copy()
$ws.Range("A38:K38").PasteSpecial(-4122)
$ws.Range("A37:K37").Copy()
$ws.Range("A39:K39").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D38").Value = 10
$ws.Range("E38").Value = "Vytvoření GUI"
$ws.Range("H38").Value = 10
$ws.Range("I38").Value = "Dopsání controllerů pro komunikaci GUI a databáze."

$ws.Range("H39").Value = 0.75
$ws.Range("I39").Value = "Sepsání uživatelského manuálu."

# I39 (unlike the rest of the new rows) keeps the plain un-bordered style.
$ws.Range("H33").Copy()
$ws.Range("I39").PasteSpecial(-4122)
$ws.Range("I39").Value = "Sepsání uživatelského manuálu."
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3) The old row 38 ("11.týden") is now row 40 - add the Petr entry that
#    was written for that week (also with the plain un-bordered style).
# ---------------------------------------------------------------------
$ws.Range("H40").Value = 3
$ws.Range("H33").Copy()
$ws.Range("I40").PasteSpecial(-4122)
$ws.Range("I40").Value = "Generování konečné upravené dokumentace a úprava výsledného dokumentu."
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 4) The "Celkem ke 4. iteraci" row is now row 41 - widen its SUM ranges
#    to cover the two new rows (and, for Petr/Lenka, to also start from
#    row 35 like the author re-keyed them).
# ---------------------------------------------------------------------
$ws.Range("B41").Formula = "=SUM(B35:B40)"
$ws.Range("D41").Formula = "=SUM(D36:D40)"
$ws.Range("F41").Formula = "=SUM(F36:F40)"
$ws.Range("H41").Formula = "=SUM(H35:H40)"
$ws.Range("J41").Formula = "=SUM(J35:J40)"

# ---------------------------------------------------------------------
# 5) Brand-new grand-total row summing all four iterations together.
# ---------------------------------------------------------------------
$ws.Range("H33").Copy()
$ws.Range("A42").PasteSpecial(-4122)
$ws.Range("D42").PasteSpecial(-4122)
$ws.Range("F42").PasteSpecial(-4122)
$ws.Range("H42").PasteSpecial(-4122)
$ws.Range("J42").PasteSpecial(-4122)
$ws.Range("B35").Copy()
$ws.Range("B42").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A42").Value = "Celkem ze všech iterací"
$ws.Range("B42").Formula = "=B41+B34+B28+D42"
$ws.Range("D42").Formula = "=D41+D34+D28+D22"
$ws.Range("F42").Formula = "=F41+F34+F28+F22"
$ws.Range("H42").Formula = "=H41+H34+H28+H22"
$ws.Range("J42").Formula = "=J41+J34+J28+J22"

# ---------------------------------------------------------------------
# 6) Update the view so it matches where the author ended up scrolled to.
# ---------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 31
$ws.Range("K42").Select()
